# Updates Sheets market-data columns (H:N) for various crafting-profit rows,
# matching the latest scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 853.3333
$ws.Range("I41").Value = 245
$ws.Range("J41").Value = 1027.1428
$ws.Range("K41").Value = 245
$ws.Range("L41").Value = 1027.1428
$ws.Range("M41").Value = 195
$ws.Range("N41").Value = -1907.1428
$ws.Range("H101").Value = 805.375
$ws.Range("J101").Value = 1195
$ws.Range("L101").Value = 3585
$ws.Range("N101").Value = -6829
$ws.Range("H107").Value = 2520
$ws.Range("I107").Value = 2900
$ws.Range("K107").Value = 2900
$ws.Range("M107").Value = -980
$ws.Range("H111").Value = 1294
$ws.Range("I111").Value = 764.5
$ws.Range("K111").Value = 2293.5
$ws.Range("M111").Value = 773.5
$ws.Range("H138").Value = 3331.5
$ws.Range("J138").Value = 3622.1633
$ws.Range("L138").Value = 10866.4899
$ws.Range("N138").Value = -21146.4899

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5515.6104
$ws.Range("I32").Value = 3182.7437
$ws.Range("J32").Value = 10064.7
$ws.Range("K32").Value = 3182.7437
$ws.Range("L32").Value = 10064.7
$ws.Range("M32").Value = -2895.7437
$ws.Range("N32").Value = -10638.7
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null
$ws.Range("H110").Value = 1204.4
$ws.Range("I110").Value = 1204.4
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1204.4
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 840.5999999999999
$ws.Range("N110").Value = $null
$ws.Range("H137").Value = 41720
$ws.Range("J137").Value = 41720
$ws.Range("L137").Value = 41720
$ws.Range("N137").Value = -51920

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 118852
$ws.Range("J59").Value = 118852
$ws.Range("L59").Value = 118852
$ws.Range("N59").Value = -120546
$ws.Range("H94").Value = 1048.6285
$ws.Range("I94").Value = 902.9259
$ws.Range("J94").Value = 1540.375
$ws.Range("K94").Value = 902.9259
$ws.Range("L94").Value = 1540.375
$ws.Range("M94").Value = -451.9259
$ws.Range("N94").Value = -2442.375
$ws.Range("H99").Value = 1578.027
$ws.Range("I99").Value = 1129.8462
$ws.Range("J99").Value = 2637.3635
$ws.Range("K99").Value = 1129.8462
$ws.Range("L99").Value = 2637.3635
$ws.Range("M99").Value = 368.1538
$ws.Range("N99").Value = -5633.363499999999
$ws.Range("H107").Value = 1562.3
$ws.Range("I107").Value = 1537.25
$ws.Range("J107").Value = 1662.5
$ws.Range("K107").Value = 1537.25
$ws.Range("L107").Value = 1662.5
$ws.Range("M107").Value = 382.75
$ws.Range("N107").Value = -5502.5
$ws.Range("H137").Value = 35453.332
$ws.Range("J137").Value = 40680
$ws.Range("L137").Value = 40680
$ws.Range("N137").Value = -50880

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22731156
$ws.Range("I31").Value = 1950.9166
$ws.Range("J31").Value = 50006200
$ws.Range("K31").Value = 1950.9166
$ws.Range("L31").Value = 50006200
$ws.Range("M31").Value = -1655.9166
$ws.Range("N31").Value = -50006790
$ws.Range("H34").Value = 22731156
$ws.Range("I34").Value = 1950.9166
$ws.Range("J34").Value = 50006200
$ws.Range("K34").Value = 1950.9166
$ws.Range("L34").Value = 50006200
$ws.Range("M34").Value = -1748.9166
$ws.Range("N34").Value = -50006604
$ws.Range("H134").Value = 4682.25
$ws.Range("I134").Value = 5160.875
$ws.Range("J134").Value = 3725
$ws.Range("K134").Value = 15482.625
$ws.Range("L134").Value = 11175
$ws.Range("M134").Value = -12947.625
$ws.Range("N134").Value = -16245

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 15462.667
$ws.Range("J39").Value = 15462.667
$ws.Range("L39").Value = 46388.001
$ws.Range("N39").Value = -46976.001
$ws.Range("H40").Value = 406.5
$ws.Range("I40").Value = 130.4
$ws.Range("J40").Value = 866.6667
$ws.Range("K40").Value = 521.6
$ws.Range("L40").Value = 3466.6668
$ws.Range("M40").Value = -452.6
$ws.Range("N40").Value = -3604.6668
$ws.Range("H46").Value = 1756
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1756
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5268
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -5450
$ws.Range("H75").Value = 2315.6667
$ws.Range("I75").Value = 914
$ws.Range("J75").Value = 2596
$ws.Range("K75").Value = 2742
$ws.Range("L75").Value = 7788
$ws.Range("M75").Value = -1744
$ws.Range("N75").Value = -9784
$ws.Range("H78").Value = 2315.6667
$ws.Range("I78").Value = 914
$ws.Range("J78").Value = 2596
$ws.Range("K78").Value = 8226
$ws.Range("L78").Value = 23364
$ws.Range("M78").Value = -3234
$ws.Range("N78").Value = -33348
$ws.Range("H109").Value = 5531.5454
$ws.Range("I109").Value = 802.8333
$ws.Range("J109").Value = 11206
$ws.Range("K109").Value = 2408.4999
$ws.Range("L109").Value = 33618
$ws.Range("M109").Value = -1368.4999
$ws.Range("N109").Value = -35698
$ws.Range("H112").Value = 5157.143
$ws.Range("I112").Value = 3266.6667
$ws.Range("J112").Value = 5672.727
$ws.Range("K112").Value = 9800.000100000001
$ws.Range("L112").Value = 17018.181
$ws.Range("M112").Value = -8692.000100000001
$ws.Range("N112").Value = -19234.181
$ws.Range("H113").Value = 558.119
$ws.Range("J113").Value = 584
$ws.Range("L113").Value = 1752
$ws.Range("N113").Value = -6092
$ws.Range("H121").Value = 1916.2115
$ws.Range("J121").Value = 1916.2115
$ws.Range("L121").Value = 5748.6345
$ws.Range("N121").Value = -8368.6345

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2266.6667
$ws.Range("I31").Value = 900
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 900
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -608
$ws.Range("N31").Value = -5584
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = $null
$ws.Range("H35").Value = 35000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 35000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 35000
$ws.Range("M35").Value = $null
$ws.Range("N35").Value = -35596
$ws.Range("H37").Value = 2266.6667
$ws.Range("I37").Value = 900
$ws.Range("J37").Value = 5000
$ws.Range("K37").Value = 900
$ws.Range("L37").Value = 5000
$ws.Range("M37").Value = -623
$ws.Range("N37").Value = -5554
$ws.Range("H41").Value = 13628.714
$ws.Range("J41").Value = 21562
$ws.Range("L41").Value = 21562
$ws.Range("N41").Value = -22272
$ws.Range("H53").Value = 26985.5
$ws.Range("J53").Value = 26985.5
$ws.Range("L53").Value = 26985.5
$ws.Range("N53").Value = -28247.5
$ws.Range("H59").Value = 8000
$ws.Range("J59").Value = 8000
$ws.Range("L59").Value = 8000
$ws.Range("N59").Value = -9166
$ws.Range("H80").Value = 25002792
$ws.Range("I80").Value = 50002464
$ws.Range("J80").Value = 3120
$ws.Range("K80").Value = 50002464
$ws.Range("L80").Value = 3120
$ws.Range("M80").Value = -50001466
$ws.Range("N80").Value = -5116
$ws.Range("H83").Value = 25002792
$ws.Range("I83").Value = 50002464
$ws.Range("J83").Value = 3120
$ws.Range("K83").Value = 250012320
$ws.Range("L83").Value = 15600
$ws.Range("M83").Value = -250007328
$ws.Range("N83").Value = -25584
$ws.Range("H137").Value = 45220
$ws.Range("J137").Value = 45220
$ws.Range("L137").Value = 45220
$ws.Range("N137").Value = -55420

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4125.6216
$ws.Range("I82").Value = 5674.75
$ws.Range("J82").Value = 2303.1177
$ws.Range("K82").Value = 5674.75
$ws.Range("L82").Value = 2303.1177
$ws.Range("M82").Value = -5313.75
$ws.Range("N82").Value = -3025.1177
$ws.Range("H85").Value = 4125.6216
$ws.Range("I85").Value = 5674.75
$ws.Range("J85").Value = 2303.1177
$ws.Range("K85").Value = 5674.75
$ws.Range("L85").Value = 2303.1177
$ws.Range("M85").Value = -4426.75
$ws.Range("N85").Value = -4799.1177

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3140.862
$ws.Range("I122").Value = 2022.7059
$ws.Range("J122").Value = 4724.9165
$ws.Range("K122").Value = 6068.1177
$ws.Range("L122").Value = 14174.7495
$ws.Range("M122").Value = -3618.1177
$ws.Range("N122").Value = -19074.7495
$ws.Range("H132").Value = 8775886
$ws.Range("I132").Value = 4650.875
$ws.Range("J132").Value = 23812288
$ws.Range("K132").Value = 13952.625
$ws.Range("L132").Value = 71436864
$ws.Range("M132").Value = -11422.625
$ws.Range("N132").Value = -71441924
